# TC06_Bento_MultiFilter_Arm-Diagnosis-TumorSize-PRStatus-EndocrineTher.xlsx
# "updated bento tc as per bento perf data availability"
#
# The four Cypher queries stored on the "startup" sheet (CasesTab, StatQuery,
# SamplesTab, FilesTab rows) all filter on tp.endocrine_therapy_type. The
# filter value is switched from "OFS" to "Tam" in every one of them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update every query cell in one shot (each cell contains exactly one
# occurrence of the endocrine_therapy_type filter value).
$ws.Cells.Replace("tp.endocrine_therapy_type IN [`"OFS`"]", "tp.endocrine_therapy_type IN [`"Tam`"]") | Out-Null

# Restore the active selection/view to D2 (matches the workbook state saved
# alongside the query-text edit).
$ws.Range("D2").Select() | Out-Null
